$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "MYSORE"
$ws.Range("C5").Value = "BANGLORE"
